# Generate Report for Handoff
# Adds a new localization-status row (for the new source file
# "0a4335b1-0133-4d88-8301-0a068dcadde5ooo....md") to all three sheets:
#   Overview, zh-cn, de-de
# and wires up the matching hyperlinks + table ranges.

$wb = $excel.ActiveWorkbook

# ---- literal strings reproduced verbatim from the target content ----
$mdFile        = "0a4335b1-0133-4d88-8301-0a068dcadde5ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$mdE2eDisplay  = "e2e\0a4335b1-0133-4d88-8301-0a068dcadde5ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$readyForHandoff = "Ready for handoff"
$dtOverviewHandoff = "2016-09-02 12:32:28"
$zhcnXlf       = "0a4335b1-0133-4d88-8301-0a068dcadde5oooooooooooooooooooooooooooooooooooooooo.60060ec60c2a819e6b7d8f3af968330953e1776e.zh-cn.xlf"
$dtZhcnHandoff = "2016-09-02 12:32:24"
$dedeXlf       = "0a4335b1-0133-4d88-8301-0a068dcadde5oooooooooooooooooooooooooooooooooooooooo.60060ec60c2a819e6b7d8f3af968330953e1776e.de-de.xlf"

$newCommit = "60060ec60c2a819e6b7d8f3af968330953e1776e"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$newCommit/e2e/$mdFile"

# =====================================================================
# Overview sheet
# =====================================================================
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A3").Value = $mdFile
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $mdUrl, "", "", $mdE2eDisplay) | Out-Null
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D3").Value = ""
$wsOv.Range("E3").Value = $readyForHandoff
$wsOv.Range("F3").Value = $readyForHandoff
$wsOv.Range("G3").Value = $dtOverviewHandoff

$wsOv.Columns.Item(5).AutoFit() | Out-Null
$wsOv.Columns.Item(6).AutoFit() | Out-Null

# =====================================================================
# zh-cn sheet
# =====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $mdFile) | Out-Null
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $zhcnXlf
$wsZh.Range("H3").Value = $dtZhcnHandoff
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = ""

$wsZh.Columns.Item(3).AutoFit() | Out-Null

# =====================================================================
# de-de sheet
# =====================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $mdFile) | Out-Null
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $dedeXlf
$wsDe.Range("H3").Value = $dtOverviewHandoff
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = ""

$wsDe.Columns.Item(3).AutoFit() | Out-Null
